# Apply the data10th.xlsx edit:
#  - Correct several mis-spelled student/father names on Sheet1.
#  - Populate a new PhotoPath (column F) for every data row with the
#    scanned photo file name "100NN.jpg" (left-aligned), including the
#    header cell F1 which already held the "PhotoPath" label.
#  - Move the current selection to C12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Corrected Name (B) / FatherName (C) values -----------------------
$names = @{
    2  = @("Umang Verma", "Mr. Satyawan Verma")
    3  = @("Roshni",      "Mr. Banke Bihari")
    4  = @("Tamanna",     "Mr. Sunil Kumar")
    5  = @("Tanvi",       "Mr. Surjeet")
    6  = @("Aarti",       "Mr. Mahadev")
    7  = @("Mahima",      "Mr. Ramesh Kumar")
    8  = @("Khushbu",     "Mr. Neeraj")
    9  = @("Deepika",     "Mr. Balvinder")
    10 = @("Anu",         "Mr. Davender")
    11 = @("Mohini",      "Mr. Sonu Soni")
    12 = @("Neelam",      "Mr. Vijay")
    13 = @("Priyanka",    "Mr. Radhey Shyam")
    14 = @("Kajal",       "Mr. Hitender")
    15 = @("Tamanna",     "Mr. Narender")
    16 = @("Anuj",        "Mr. Parmod")
    17 = @("Manit",       "Mr. Karan Singh")
    18 = @("Nitin",       "Mr. Mangat Ram")
    19 = @("Ranveer",     "Mr. Amit")
    20 = @("Dev",         "Mr. Mahender")
    21 = @("Vansh",       "Mr. Dhanajay")
    22 = @("Pushp",       "Mr. Nehru Dutt")
    23 = @("Jai Soni",    "Mr. Ravinder")
    24 = @("Dinesh",      "Mr. Jai Singh")
    25 = @("Sanyam",      "Mr. Anil Kumar")
    26 = @("Mohit",       "Mr. Harish ")
    27 = @("Manav",       "Mr. Mohit")
    28 = @("Prince",      "Mr. Rajesh")
    29 = @("Dheeraj",     "Mr. Rinku")
    30 = @("Sanjana",     "Mr. Satish Kumar")
}

foreach ($row in $names.Keys) {
    $pair = $names[$row]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}

# --- New PhotoPath column (F), left-aligned, rows 1 (header) to 30 ----
$ws.Range("F1:F30").HorizontalAlignment = -4131   # xlLeft

for ($row = 2; $row -le 30; $row++) {
    $photoId = 10000 + ($row - 1)
    $ws.Cells.Item($row, 6).Value = "$photoId.jpg"
}

# --- Selection -----------------------------------------------------
$ws.Range("C12").Select() | Out-Null
